# Update "想去人数" (F column) counts on the 展览, 演出 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 15023
$ws1.Range("F3").Value  = 19081
$ws1.Range("F14").Value = 163
$ws1.Range("F15").Value = 221
$ws1.Range("F17").Value = 1472
$ws1.Range("F21").Value = 237
$ws1.Range("F22").Value = 7956
$ws1.Range("F32").Value = 170
$ws1.Range("F34").Value = 287
$ws1.Range("F35").Value = 5440
$ws1.Range("F36").Value = 469
$ws1.Range("F37").Value = 13
$ws1.Range("F38").Value = 28

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 16

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 15023
$ws4.Range("F3").Value  = 19081
$ws4.Range("F14").Value = 163
$ws4.Range("F15").Value = 221
$ws4.Range("F17").Value = 1472
$ws4.Range("F22").Value = 237
$ws4.Range("F23").Value = 7956
$ws4.Range("F30").Value = 16
$ws4.Range("F35").Value = 170
$ws4.Range("F37").Value = 287
$ws4.Range("F38").Value = 5440
$ws4.Range("F39").Value = 469
$ws4.Range("F40").Value = 13
$ws4.Range("F41").Value = 28
